# Daily attendance processing - 2025-11-29 05:50:41
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) for rows where System is listed first
# and exactly one other recorder (email) is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($value -eq "System, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, System"
    }
}
